# Actualizaciones al plan del proyecto.
#
# The "schedule" sheet pulled its task details (name, criteria, estimated
# hours/week and role hour allocations) from an external workbook
# (tspi/plan.xlsx) via LOOKUP formulas. This converts those formulas to
# their last-calculated static values and removes the now-unused external
# workbook link.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("schedule")

# Columns B through L (2 through 12) on rows 2-23 contain the
# LOOKUP(...) formulas referencing the external workbook.
for ($r = 2; $r -le 23; $r++) {
    for ($col = 2; $col -le 12; $col++) {
        $cell = $ws.Cells.Item($r, $col)
        $f = $cell.Formula
        if ($f -ne $null -and $f.ToString().StartsWith("=")) {
            $cell.Value = $cell.Value2
        }
    }
}

# With every formula referencing the external workbook replaced by a
# static value, break the (now unused) external link so it is removed
# from the workbook entirely.
$links = $wb.LinkSources(1)
if ($links -ne $null) {
    foreach ($link in $links) {
        $wb.BreakLink($link, 1)
    }
}
